# All Result ready to start write
# Remove the old column A (row-index values 11/15, previously header-styled
# but without a header label) and shift B:F left into A:E. The former
# column B header ("QS_Astral15") lands in A1, completing the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()
